{"js": "// Adds the 10 new \"List Paragraph\" bullet items (numbered list numId=1) that\n// document the Board/BoardScanner/BoardManager/Renderer refactor, right\n// after the existing \"Menu() method renamed to Run();\" bullet and before\n// the section break at the end of the document body.\n\nconst newBullets = [\n  \"In Board.cs fields renamed to fieldsMatrix (to be more descriptive)\",\n  \"Add IBoard interface (Board implements it)\",\n  \"IBoardManager interface created\",\n  \"IBoardScanner interface created (there is repeating logic, so it will be extracted in separate class implementing this interface)\",\n  \"BoardScanner.cs created.\",\n  \"Replace logic for ScanSurroundingFields from Board to the new BoardScanner class. It depends on the private IsMineInPosition() method, so it goes to the new class too.\",\n  \"BoardManager.cs created implemented IBoardManager \",\n  \"Renderer.cs created implementing IRenderer\",\n  \"Duplicated methods wich are extracted in separate classes are removed from Board.cs\",\n  \"Adapt Game.cs work with new classes\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the body is \"Menu() method renamed to Run();\" --\n// insert all new bullet paragraphs right after it, one by one, so each\n// becomes the new \"anchor\" for the next insertion (keeps them in order).\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\nfor (const text of newBullets) {\n  const newPara = anchor.insertParagraph(text, Word.InsertLocation.after);\n  // Match the existing numbered-list bullet formatting used throughout the\n  // \"What's new\" log (pStyle \"a3\" == built-in \"List Paragraph\" style, with\n  // the same list -- numId 1 / ilvl 0 -- as the preceding bullets).\n  newPara.style = \"List Paragraph\";\n  newPara.attachToList(1, 0);\n  anchor = newPara;\n}\n\nawait context.sync();\n", "ps1": "# Adds the 10 new \"List Paragraph\" bullet items (numbered list numId=1) that\n# document the Board/BoardScanner/BoardManager/Renderer refactor, right\n# after the existing \"Menu() method renamed to Run();\" bullet and before\n# the section break at the end of the document body.\n\n$d = $word.ActiveDocument\n\n$newBullets = @(\n  \"In Board.cs fields renamed to fieldsMatrix (to be more descriptive)\",\n  \"Add IBoard interface (Board implements it)\",\n  \"IBoardManager interface created\",\n  \"IBoardScanner interface created (there is repeating logic, so it will be extracted in separate class implementing this interface)\",\n  \"BoardScanner.cs created.\",\n  \"Replace logic for ScanSurroundingFields from Board to the new BoardScanner class. It depends on the private IsMineInPosition() method, so it goes to the new class too.\",\n  \"BoardManager.cs created implemented IBoardManager \",\n  \"Renderer.cs created implementing IRenderer\",\n  \"Duplicated methods wich are extracted in separate classes are removed from Board.cs\",\n  \"Adapt Game.cs work with new classes\"\n)\n\n# Reuse the list template already applied to the existing numbered bullets\n# (numId 1) so the new paragraphs continue the very same list.\n$existingListPara = $d.Paragraphs.Item(2)\n$listTemplate = $existingListPara.Range.ListFormat.ListTemplate\n\nforeach ($bulletText in $newBullets) {\n    $lastPara = $d.Paragraphs.Last\n    $r = $lastPara.Range\n    $r.Collapse(0)\n    $r.InsertAfter([char]13 + $bulletText)\n\n    $newPara = $d.Paragraphs.Last\n    $newPara.Style = \"List Paragraph\"\n    $newPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 0, $false, $false)\n}\n"}
